$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI TPM results: data rows 2-13 recalculated, and 3 new rows (14-16)
# added for the new "Resolving-Mac" target cluster (one per sending cluster).
$rows = @(
  @("ECs", "Agt", "Mas1", "ECs", 2, [double]"0.6666666666666666", [double]"0.07702833333333334", [double]"0.231085", [double]"0.02259036512642383", [double]"0.02259036512642383", 2, [double]"0.6666666666666666", [double]"0.4865679999999999", [double]"1.459704", [double]"0.2260283549680654", [double]"0.2260283549680655", [double]"0.03747952209333333", [double]"0.33731569884", [double]"0.005106063067653532", [double]"0.005106063067653534"),
  @("ECs", "Agt", "Mas1", "FAPs", 2, [double]"0.6666666666666666", [double]"0.07702833333333334", [double]"0.231085", [double]"0.02259036512642383", [double]"0.02259036512642383", 3, 1, [double]"1.254922333333333", [double]"3.764767", [double]"0.5829566075369108", [double]"0.582956607536911", [double]"0.09666457579944444", [double]"0.8699811821950001", [double]"0.01316920261712017", [double]"0.01316920261712018"),
  @("ECs", "Agt", "Mas1", "Inflammatory-Mac", 2, [double]"0.6666666666666666", [double]"0.07702833333333334", [double]"0.231085", [double]"0.02259036512642383", [double]"0.02259036512642383", 1, [double]"0.3333333333333333", [double]"0.000891", [double]"0.002673", [double]"0.0004139015806147266", [double]"0.0004139015806147268", [double]"6.8632245E-05", [double]"0.000617690205", [double]"9.350187832490623E-06", [double]"9.350187832490627E-06"),
  @("ECs", "Agt", "Mas1", "MuSCs", 2, [double]"0.6666666666666666", [double]"0.07702833333333334", [double]"0.231085", [double]"0.02259036512642383", [double]"0.02259036512642383", 3, 1, [double]"0.4094323333333333", [double]"1.228297", [double]"0.1901960605178926", [double]"0.1901960605178926", [double]"0.03153789024944445", [double]"0.283841012245", [double]"0.004296598452706597", [double]"0.004296598452706598"),
  @("ECs", "Agt", "Mas1", "Resolving-Mac", 2, [double]"0.6666666666666666", [double]"0.07702833333333334", [double]"0.231085", [double]"0.02259036512642383", [double]"0.02259036512642383", 1, [double]"0.3333333333333333", [double]"0.0008719999999999999", [double]"0.002616", [double]"0.0004050753965163206", [double]"0.0004050753965163207", [double]"6.716870666666666E-05", [double]"0.00060451836", [double]"9.150801111034595E-06", [double]"9.150801111034596E-06"),
  @("FAPs", "Agt", "Mas1", "ECs", 3, 1, [double]"3.276419666666667", [double]"9.829259", [double]"0.9608869019286738", [double]"0.9608869019286738", 2, [double]"0.6666666666666666", [double]"0.4865679999999999", [double]"1.459704", [double]"0.2260283549680654", [double]"0.2260283549680655", [double]"1.594200964370667", [double]"14.347808679336", [double]"0.217187685753299", [double]"0.217187685753299"),
  @("FAPs", "Agt", "Mas1", "FAPs", 3, 1, [double]"3.276419666666667", [double]"9.829259", [double]"0.9608869019286738", [double]"0.9608869019286738", 3, 1, [double]"1.254922333333333", [double]"3.764767", [double]"0.5829566075369108", [double]"0.582956607536911", [double]"4.111652213072555", [double]"37.004869917653", [double]"0.5601553685749919", [double]"0.5601553685749922"),
  @("FAPs", "Agt", "Mas1", "Inflammatory-Mac", 3, 1, [double]"3.276419666666667", [double]"9.829259", [double]"0.9608869019286738", [double]"0.9608869019286738", 1, [double]"0.3333333333333333", [double]"0.000891", [double]"0.002673", [double]"0.0004139015806147266", [double]"0.0004139015806147268", [double]"0.002919289923", [double]"0.026273609307", [double]"0.0003977126075002659", [double]"0.0003977126075002661"),
  @("FAPs", "Agt", "Mas1", "MuSCs", 3, 1, [double]"3.276419666666667", [double]"9.829259", [double]"0.9608869019286738", [double]"0.9608869019286738", 3, 1, [double]"0.4094323333333333", [double]"1.228297", [double]"0.1901960605178926", [double]"0.1901960605178926", [double]"1.341472149102556", [double]"12.073249341923", [double]"0.1827569033500764", [double]"0.1827569033500764"),
  @("FAPs", "Agt", "Mas1", "Resolving-Mac", 3, 1, [double]"3.276419666666667", [double]"9.829259", [double]"0.9608869019286738", [double]"0.9608869019286738", 1, [double]"0.3333333333333333", [double]"0.0008719999999999999", [double]"0.002616", [double]"0.0004050753965163206", [double]"0.0004050753965163207", [double]"0.002857037949333334", [double]"0.025713341544", [double]"0.0003892316428060964", [double]"0.0003892316428060965"),
  @("MuSCs", "Agt", "Mas1", "ECs", 1, [double]"0.3333333333333333", [double]"0.056339", [double]"0.169017", [double]"0.01652273294490242", [double]"0.01652273294490242", 2, [double]"0.6666666666666666", [double]"0.4865679999999999", [double]"1.459704", [double]"0.2260283549680654", [double]"0.2260283549680655", [double]"0.027412754552", [double]"0.246714790968", [double]"0.003734606147112954", [double]"0.003734606147112955"),
  @("MuSCs", "Agt", "Mas1", "FAPs", 1, [double]"0.3333333333333333", [double]"0.056339", [double]"0.169017", [double]"0.01652273294490242", [double]"0.01652273294490242", 3, 1, [double]"1.254922333333333", [double]"3.764767", [double]"0.5829566075369108", [double]"0.582956607536911", [double]"0.07070106933766666", [double]"0.636309624039", [double]"0.009632036344798667", [double]"0.009632036344798671"),
  @("MuSCs", "Agt", "Mas1", "Inflammatory-Mac", 1, [double]"0.3333333333333333", [double]"0.056339", [double]"0.169017", [double]"0.01652273294490242", [double]"0.01652273294490242", 1, [double]"0.3333333333333333", [double]"0.000891", [double]"0.002673", [double]"0.0004139015806147266", [double]"0.0004139015806147268", [double]"5.0198049E-05", [double]"0.000451782441", [double]"6.83878528197013E-06", [double]"6.838785281970132E-06"),
  @("MuSCs", "Agt", "Mas1", "MuSCs", 1, [double]"0.3333333333333333", [double]"0.056339", [double]"0.169017", [double]"0.01652273294490242", [double]"0.01652273294490242", 3, 1, [double]"0.4094323333333333", [double]"1.228297", [double]"0.1901960605178926", [double]"0.1901960605178926", [double]"0.02306700822766667", [double]"0.207603074049", [double]"0.003142558715109639", [double]"0.00314255871510964"),
  @("MuSCs", "Agt", "Mas1", "Resolving-Mac", 1, [double]"0.3333333333333333", [double]"0.056339", [double]"0.169017", [double]"0.01652273294490242", [double]"0.01652273294490242", 1, [double]"0.3333333333333333", [double]"0.0008719999999999999", [double]"0.002616", [double]"0.0004050753965163206", [double]"0.0004050753965163207", [double]"4.9127608E-05", [double]"0.000442148472", [double]"6.692952599189623E-06", [double]"6.692952599189624E-06")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
  $rowNum = $i + 2
  $vals = $rows[$i]
  for ($j = 0; $j -lt $vals.Count; $j++) {
    $ws.Cells.Item($rowNum, $j + 1).Value = $vals[$j]
  }
}
